$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default/general "Normal" style, used to restore the
# style index on cells we have to briefly force to Text format (below) so their
# literal digit-string content (e.g. "321.07") is not reinterpreted as a number.
$normalStyle = $ws.Range("C2").Style

$ws.Range("D2").Value = "27.735.79"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.861.16"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  -1.03%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "321.07"
$cell.Style = $normalStyle
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  -0.95%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4373"
$cell.Style = $normalStyle
$ws.Range("E7").Value = "  -0.47%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3782"
$cell.Style = $normalStyle
$ws.Range("E8").Value = "  -0.12%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07416"
$cell.Style = $normalStyle
$ws.Range("E9").Value = "  +0.06%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.8853"
$cell.Style = $normalStyle
$ws.Range("E10").Value = "  +1.06%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "21.63"
$cell.Style = $normalStyle
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.859.30"
$ws.Range("E12").Value = "  +0.19%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "6.759"
$cell.Style = $normalStyle
$ws.Range("E13").Value = "  +0.94%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.502"
$cell.Style = $normalStyle
$ws.Range("E14").Value = "  -0.33%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.07148"
$cell.Style = $normalStyle
$ws.Range("E15").Value = "  -0.92%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "87.75"
$cell.Style = $normalStyle
$ws.Range("E16").Value = "  +5.49%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "1.024"
$cell.Style = $normalStyle
$ws.Range("E17").Value = "  -1.08%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000009054"
$cell.Style = $normalStyle
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  -0.93%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "15.49"
$cell.Style = $normalStyle
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "27.768.75"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("E22").Value = "  +0.65%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "11.17"
$cell.Style = $normalStyle
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("D24").Value = "2.094.11"
$ws.Range("E24").Value = "  +0.95%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.039"
$cell.Style = $normalStyle
$ws.Range("E25").Value = "  +6.18%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "157.26"
$cell.Style = $normalStyle
$ws.Range("E26").Value = "  -0.45%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "18.71"
$cell.Style = $normalStyle
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "5.431"
$cell.Style = $normalStyle
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("E29").Value = "  +0.92%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "121.62"
$cell.Style = $normalStyle
$ws.Range("E30").Value = "  +3.86%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.09067"
$cell.Style = $normalStyle
$ws.Range("E31").Value = "  +0.10%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.218"
$cell.Style = $normalStyle
$ws.Range("E32").Value = "  +1.61%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.7705"
$cell.Style = $normalStyle
$ws.Range("E33").Value = "  +1.30%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.037"
$cell.Style = $normalStyle
$ws.Range("E34").Value = "  +5.38%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.562"
$cell.Style = $normalStyle
$ws.Range("E35").Value = "  +0.82%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.021"
$cell.Style = $normalStyle
$ws.Range("E36").Value = "  -0.86%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.139"
$cell.Style = $normalStyle
$ws.Range("E37").Value = "  -0.90%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01979"
$cell.Style = $normalStyle
$ws.Range("E38").Value = "  +0.17%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.05312"
$cell.Style = $normalStyle
$ws.Range("E39").Value = "  +0.11%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.869"
$cell.Style = $normalStyle
$ws.Range("E40").Value = "  +1.92%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.5186"
$cell.Style = $normalStyle
$ws.Range("E41").Value = "  +0.65%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "6.985"
$cell.Style = $normalStyle
$ws.Range("E42").Value = "  +3.37%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.1679"
$cell.Style = $normalStyle
$ws.Range("E43").Value = "  +0.19%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "8.724"
$cell.Style = $normalStyle
$ws.Range("E44").Value = "  +2.91%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "10.80"
$cell.Style = $normalStyle
$ws.Range("E45").Value = "  +2.06%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "110.27"
$cell.Style = $normalStyle
$ws.Range("E46").Value = "  +1.42%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.714"
$cell.Style = $normalStyle
$ws.Range("E47").Value = "  +0.23%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.4732"
$cell.Style = $normalStyle
$ws.Range("E48").Value = "  +1.81%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.021"
$cell.Style = $normalStyle
$ws.Range("E49").Value = "  -0.98%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.06482"
$cell.Style = $normalStyle
$ws.Range("E50").Value = "  +1.20%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.851"
$cell.Style = $normalStyle
$ws.Range("E51").Value = "  +0.21%  "
